# Auto-generated edit script for cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) cells must be forced to Text format before assignment ---
# so that numeric-looking strings (e.g. "1.000", "28.423.67") are stored verbatim
# as text, matching the source data which is text, not numeric.
$priceCells = @("D2", "D3", "D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D45", "D50", "D51")
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "28.423.67"
$ws.Range("E2").Value = "  +0.66%  "

# Row 3
$ws.Range("D3").Value = "1.829.98"
$ws.Range("E3").Value = "  +2.37%  "

# Row 4
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").Value = "317.96"
$ws.Range("E5").Value = "  +0.56%  "

# Row 6
$ws.Range("E6").Value = "  -0.10%  "

# Row 7
$ws.Range("D7").Value = "0.5318"
$ws.Range("E7").Value = "  +0.11%  "

# Row 8
$ws.Range("D8").Value = "0.4055"
$ws.Range("E8").Value = "  +8.06%  "

# Row 9
$ws.Range("D9").Value = "0.07627"
$ws.Range("E9").Value = "  +1.97%  "

# Row 10
$ws.Range("D10").Value = "41.83"
$ws.Range("E10").Value = "  +0.81%  "

# Row 11
$ws.Range("D11").Value = "1.108"
$ws.Range("E11").Value = "  +1.47%  "

# Row 12
$ws.Range("D12").Value = "6.358"
$ws.Range("E12").Value = "  +4.36%  "

# Row 13
$ws.Range("B13").Value = "BinanceUSD"
$ws.Range("C13").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D13").Value = "1.001"
$ws.Range("E13").Value = "  -0.14%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "7.617"
$ws.Range("E14").Value = "  +5.27%  "

# Row 15
$ws.Range("B15").Value = "Solana"
$ws.Range("C15").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D15").Value = "20.92"
$ws.Range("E15").Value = "  +2.45%  "

# Row 16
$ws.Range("D16").Value = "1.828.01"
$ws.Range("E16").Value = "  +1.43%  "

# Row 17
$ws.Range("D17").Value = "89.41"
$ws.Range("E17").Value = "  +0.31%  "

# Row 18
$ws.Range("D18").Value = "0.00001075"
$ws.Range("E18").Value = "  +1.85%  "

# Row 19
$ws.Range("D19").Value = "0.06616"
$ws.Range("E19").Value = "  +1.88%  "

# Row 20
$ws.Range("E20").Value = "  +1.32%  "

# Row 21
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.11%  "

# Row 22
$ws.Range("D22").Value = "6.098"
$ws.Range("E22").Value = "  +2.95%  "

# Row 23
$ws.Range("D23").Value = "28.444.42"
$ws.Range("E23").Value = "  +0.58%  "

# Row 24
$ws.Range("D24").Value = "11.24"
$ws.Range("E24").Value = "  +1.14%  "

# Row 25
$ws.Range("D25").Value = "2.147"
$ws.Range("E25").Value = "  +2.53%  "

# Row 26
$ws.Range("D26").Value = "2.482"
$ws.Range("E26").Value = "  +8.49%  "

# Row 27
$ws.Range("D27").Value = "157.29"
$ws.Range("E27").Value = "  -0.54%  "

# Row 28
$ws.Range("D28").Value = "20.58"
$ws.Range("E28").Value = "  +1.53%  "

# Row 29
$ws.Range("D29").Value = "2.039.74"
$ws.Range("E29").Value = "  +1.97%  "

# Row 30
$ws.Range("D30").Value = "124.49"
$ws.Range("E30").Value = "  +2.83%  "

# Row 31
$ws.Range("D31").Value = "1.128"
$ws.Range("E31").Value = "  +2.91%  "

# Row 32
$ws.Range("D32").Value = "0.1093"
$ws.Range("E32").Value = "  +4.85%  "

# Row 33
$ws.Range("D33").Value = "5.683"
$ws.Range("E33").Value = "  +3.06%  "

# Row 34
$ws.Range("D34").Value = "3.655"
$ws.Range("E34").Value = "  -0.15%  "

# Row 35
$ws.Range("D35").Value = "0.07149"
$ws.Range("E35").Value = "  +11.86%  "

# Row 36
$ws.Range("D36").Value = "0.2263"
$ws.Range("E36").Value = "  +0.94%  "

# Row 37
$ws.Range("D37").Value = "0.02343"
$ws.Range("E37").Value = "  +3.02%  "

# Row 38
$ws.Range("D38").Value = "5.218"
$ws.Range("E38").Value = "  +4.46%  "

# Row 39
$ws.Range("D39").Value = "8.843"
$ws.Range("E39").Value = "  +4.25%  "

# Row 40
$ws.Range("D40").Value = "0.6279"
$ws.Range("E40").Value = "  +1.85%  "

# Row 41
$ws.Range("E41").Value = "  +2.70%  "

# Row 42
$ws.Range("D42").Value = "1.187"
$ws.Range("E42").Value = "  +1.08%  "

# Row 43
$ws.Range("E43").Value = "  -0.09%  "

# Row 44
$ws.Range("E44").Value = "  -2.07%  "

# Row 45
$ws.Range("D45").Value = "13.49"
$ws.Range("E45").Value = "  +1.55%  "

# Row 46
$ws.Range("E46").Value = "  +1.04%  "

# Row 47
$ws.Range("E47").Value = "  +1.63%  "

# Row 48
$ws.Range("E48").Value = "  +0.06%  "

# Row 49
$ws.Range("E49").Value = "  +3.10%  "

# Row 50
$ws.Range("D50").Value = "1.201"
$ws.Range("E50").Value = "  +0.04%  "

# Row 51
$ws.Range("D51").Value = "0.06899"
$ws.Range("E51").Value = "  +0.85%  "
